# Apply updated crypto price/volume data to columns D (Price) and E (Volume 1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.164.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.449.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.37%  "

$ws.Range("E7").Value = "  +1.23%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +4.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.83%  "

$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.123"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.86%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.41%  "

$ws.Range("E14").Value = "  +1.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.831.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.454.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.16%  "

$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.019.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0932"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "246.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.33%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("E27").Value = "  +2.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.38%  "

$ws.Range("E34").Value = "  +3.22%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  -1.43%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "127.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.48%  "

$ws.Range("E41").Value = "  +1.73%  "

$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.87%  "

$ws.Range("E44").Value = "  +1.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.966.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.74%  "

$ws.Range("E46").Value = "  +0.72%  "

$ws.Range("E47").Value = "  -4.83%  "

$ws.Range("E48").Value = "  +13.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.96%  "
